# Sprint 3 artifacts: populate "Sprint No." (column F) for the Sprint-3
# requirement rows (15-22) with the story's sprint number (3), matching the
# already-populated Sprint-1 (rows 3-8, value 1) and Sprint-2 (rows 9-14,
# value 2) rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 15; $row -le 22; $row++) {
    $ws.Cells.Item($row, 6).Value2 = 3
}
